$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: merge "A" + " " + "slide" runs into a single run "A slide"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "A slide"

# Table cell (second cell of first row) : merge "a" + " " + "table" runs into a single run "a table"
$tbl = $s.Shapes.Item(3).Table
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "a table"
